# Update the "winning team" picks on the "Round 1" sheet of the March
# Madness bracket workbook. Column B holds each player's predicted winner
# for a given matchup; it was previously blank for ten of the rows. Column
# A (losing team, via formula), and the per-player "correct pick" flags in
# D/F/H/J/L (also formulas) recompute automatically once B is filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Round 1")

$ws.Range("B10").Value = "Gonz"
$ws.Range("B11").Value = "Baylor"
$ws.Range("B12").Value = "Murray St"
$ws.Range("B16").Value = "FL"
$ws.Range("B17").Value = "Michigan"
$ws.Range("B22").Value = "Nova"
$ws.Range("B23").Value = "Purdue"
$ws.Range("B29").Value = "Kansas"
$ws.Range("B32").Value = "Wofford"
$ws.Range("B33").Value = "Kentucky"

# New conditional format: highlight a pick in column B when it matches the
# value typed in U1 (a quick visual "find my pick" helper). Inserted with
# top priority, same as Excel does for a freshly added rule.
$rng = $ws.Range("B2:B33")
$fc = $rng.FormatConditions.Add(1, 3, "=`$U`$1")
$fc.SetFirstPriority()
$fc.Interior.Color = 12566463
